$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for D/E columns so numeric-looking strings
# (prices, percentages) are preserved as text, matching the source data.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "309.32"
$ws.Range("E2").Value = "-2.79%"
$ws.Range("D3").Value = "37.54"
$ws.Range("E3").Value = "-5.68%"
$ws.Range("D4").Value = "5.125"
$ws.Range("E4").Value = "-0.34%"
$ws.Range("D5").Value = "0.07861"
$ws.Range("E5").Value = "-4.27%"
$ws.Range("D6").Value = "1.961"
$ws.Range("E6").Value = "-8.11%"
$ws.Range("D7").Value = "4.364"
$ws.Range("E7").Value = "1.84%"
$ws.Range("D8").Value = "8.298"
$ws.Range("E8").Value = "-0.07%"
$ws.Range("D9").Value = "3.149"
$ws.Range("E9").Value = "-5.11%"
$ws.Range("D10").Value = "0.9266"
$ws.Range("E10").Value = "-0.89%"
$ws.Range("E11").Value = "-2.82%"
$ws.Range("D12").Value = "0.1976"
$ws.Range("E12").Value = "-0.51%"
$ws.Range("D13").Value = "0.08984"
$ws.Range("E13").Value = "-1.33%"
$ws.Range("D14").Value = "0.03444"
$ws.Range("E14").Value = "-1.02%"
$ws.Range("D15").Value = "0.09700"
$ws.Range("E15").Value = "-1.01%"
$ws.Range("E16").Value = "-0.84%"
$ws.Range("D17").Value = "0.006093"
$ws.Range("E17").Value = "0.26%"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").Value = "3.609"
$ws.Range("E18").Value = "-1.80%"
$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D19").Value = "0.3465"
$ws.Range("E19").Value = "-0.26%"
$ws.Range("B20").Value = "ProBitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D20").Value = "0.1294"
$ws.Range("E20").Value = "0.25%"
$ws.Range("B21").Value = "MCDex"
$ws.Range("C21").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D21").Value = "5.010"
$ws.Range("E21").Value = "2.29%"
$ws.Range("B22").Value = "ZBToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D22").Value = "0.2513"
$ws.Range("E22").Value = "2.61%"
$ws.Range("B23").Value = "UpBots"
$ws.Range("C23").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D23").Value = "0.02108"
$ws.Range("E23").Value = "5,174.32%"
$ws.Range("D24").Value = "0.04349"
$ws.Range("E24").Value = "0.62%"
$ws.Range("D25").Value = "0.001222"
$ws.Range("E25").Value = "-0.36%"
$ws.Range("D26").Value = "0.004542"
$ws.Range("E26").Value = "-4.78%"
$ws.Range("D27").Value = "0.0001353"
$ws.Range("E27").Value = "4.13%"
$ws.Range("D39").Value = "0.02285"
$ws.Range("E39").Value = "2.26%"
$ws.Range("D40").Value = "0.05046"
$ws.Range("E40").Value = "-3.42%"
$ws.Range("D41").Value = "0.007617"
$ws.Range("E41").Value = "1.62%"
$ws.Range("D42").Value = "0.009134"
$ws.Range("E42").Value = "-5.54%"
$ws.Range("D43").Value = "0.1356"
$ws.Range("E43").Value = "-2.15%"
$ws.Range("D44").Value = "0.002054"
$ws.Range("E44").Value = "-4.39%"
$ws.Range("D45").Value = "0.008391"
$ws.Range("E45").Value = "-11.76%"
$ws.Range("D46").Value = "0.00006774"
$ws.Range("E46").Value = "2.49%"
$ws.Range("D47").Value = "0.00000000752"
$ws.Range("E47").Value = "0.27%"
$ws.Range("D48").Value = "0.003006"
$ws.Range("E48").Value = "8.51%"
$ws.Range("D50").Value = "0.00002104"
$ws.Range("E50").Value = "0.27%"
$ws.Range("D51").Value = "0.0002004"
$ws.Range("E51").Value = "0.27%"
